$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG_TO ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-1650291288582433"
$ws1.Range("B2").Value = "go_stims-16502912885284307.csv"
$ws1.Range("B3").Value = "GNG_stims-1650291288549436.csv"
$ws1.Range("B4").Value = "go_stims-16502912885514326.csv"
$ws1.Range("B5").Value = "GNG_stims-16502912885804398.csv"

# --- Sheet 2: NB_TO ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-16502912907038465"
$ws2.Range("B2").Value = "OB-16502912894987152.csv"
$ws2.Range("B3").Value = "ZB-match_2-16502912886774335.csv"
$ws2.Range("B4").Value = "OB-16502912890547054.csv"
$ws2.Range("B5").Value = "TB-16502912906338472.csv"
$ws2.Range("B6").Value = "TB-16502912901048453.csv"
$ws2.Range("B7").Value = "OB-16502912889397085.csv"
$ws2.Range("B8").Value = "ZB-match_9-16502912886434333.csv"
$ws2.Range("B9").Value = "TB-16502912906818433.csv"
$ws2.Range("B10").Value = "ZB-match_7-16502912887074296.csv"

# --- Sheet 3: RS_TO ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-16502912907058432"
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# --- Sheet 4: TOL_TO ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-1650291290750841"
$ws4.Range("B2").Value = "MM_stims-1650291290719845.csv"
$ws4.Range("B3").Value = "ZM_stims-1650291290708842.csv"
$ws4.Range("B4").Value = "MM_stims-1650291290734841.csv"
$ws4.Range("B5").Value = "ZM_stims-16502912907208426.csv"
$ws4.Range("B6").Value = "MM_stims-16502912907498412.csv"
$ws4.Range("B7").Value = "ZM_stims-16502912907358422.csv"

# --- Sheet 5: vSAT_TO ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16502912908288417"
$ws5.Range("B2").Value = "SAT_stims-16502912907548425.csv"
$ws5.Range("B3").Value = "vSAT_stims-16502912907978404.csv"
$ws5.Range("B4").Value = "SAT_stims-16502912907818418.csv"
$ws5.Range("B5").Value = "vSAT_stims-16502912908148568.csv"
